$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.371296882629395
$ws.Range("B1").Value = 1.472806692123413
$ws.Range("C1").Value = 5.077484130859375
$ws.Range("D1").Value = 2.784038782119751
$ws.Range("E1").Value = 0.9527737498283386
